$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (keeping cell contents/positions otherwise unchanged)
$ws.Range("B1").Value = "GDP_GROWTH"
$ws.Range("C1").Value = "EXPORT_INDEX"
$ws.Range("D1").Value = "UNEMPLOYMENT"

# Row 1 height shrinks now that the headers are shorter (no longer wraps to 5 lines)
$ws.Rows.Item(1).RowHeight = 33

# Move the active selection to D2
$ws.Range("D2").Select()
